$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2990
$ws1.Range("F3").Value = 6424
$ws1.Range("F6").Value = 545
$ws1.Range("F7").Value = 82
$ws1.Range("F9").Value = 2941
$ws1.Range("F12").Value = 7614
$ws1.Range("F13").Value = 370
$ws1.Range("F18").Value = 9
$ws1.Range("F20").Value = 9311
$ws1.Range("F27").Value = 127
$ws1.Range("F28").Value = 122
$ws1.Range("F30").Value = 125
$ws1.Range("F36").Value = 2050
$ws1.Range("F37").Value = 1489
$ws1.Range("F38").Value = 783
$ws1.Range("F39").Value = 3949
$ws1.Range("F40").Value = 217
$ws1.Range("F41").Value = 46
$ws1.Range("F43").Value = 103
$ws1.Range("F45").Value = 39
$ws1.Range("F46").Value = 14
$ws1.Range("F47").Value = 66
$ws1.Range("F48").Value = 41
$ws1.Range("F49").Value = 62

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 151
$ws2.Range("G7").Value = 280
$ws2.Range("F15").Value = 7
$ws2.Range("F19").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2990
$ws4.Range("F6").Value = 6424
$ws4.Range("F8").Value = 151
$ws4.Range("G8").Value = 280
$ws4.Range("F10").Value = 545
$ws4.Range("F11").Value = 82
$ws4.Range("F13").Value = 2941
$ws4.Range("F17").Value = 7614
$ws4.Range("F18").Value = 370
$ws4.Range("F22").Value = 9
$ws4.Range("F24").Value = 9311
$ws4.Range("F29").Value = 127
$ws4.Range("F30").Value = 122
$ws4.Range("F32").Value = 125
$ws4.Range("F36").Value = 2050
$ws4.Range("F37").Value = 1489
$ws4.Range("F38").Value = 783
$ws4.Range("F40").Value = 3949
$ws4.Range("F41").Value = 217
$ws4.Range("F42").Value = 46
$ws4.Range("F44").Value = 103
$ws4.Range("F46").Value = 39
$ws4.Range("F47").Value = 66
$ws4.Range("F48").Value = 41
$ws4.Range("F49").Value = 62
